$d = $word.ActiveDocument

# --- Add the three new character styles used by the edited paragraphs ---
$styles = $d.Styles

$gaNStyle = $styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Perioadele campaniei ..." run (4 occurrences) ---
$campaignText = "Perioadele campaniei din 2022 pentru constelația Orion: 16-25 ianuarie, 14-23 februarie, 14-24 martie"
$rng = $d.Content
$found = $rng.Find.Execute($campaignText)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($campaignText)
}

# --- Apply GaNParagraph to the activity-description paragraph run ---
$paragraphText = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  constelația Orion pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($paragraphText)
if ($found2) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the credit/link run ---
$linksText = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($linksText)
if ($found3) {
    $rng3.Style = "GaNLinks"
}
